# Applies the "Updated BSD - Protocol" edit:
#  - Adds a new task ("Plakat für TDOT 2016/2017 machen") on 2017-01-25
#    14:10-15:50, status 100%, into all three "TDOT" protocol tables that
#    still had an empty row at that date slot:
#       * Tabelle245 (W6:AA40)  -> row 25
#       * Tabelle2   (C6:G39)   -> row 30
#       * Tabelle24  (M6:Q40)   -> row 30
#  - Updates the sheet view's top-left cell / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taskName = "Plakat für TDOT 2016/2017 machen"
$taskDate = 42760                    # 2017-01-25 (serial date)
$taskFrom = 0.59027777777777779      # 14:10
$taskTo   = 0.65972222222222221      # 15:50
$taskStatus = 1

# --- Row 25: third table (W:AA) ---
$ws.Range("W25").Value = $taskName
$ws.Range("X25").Value = $taskDate
$ws.Range("Y25").Value = $taskFrom
$ws.Range("Z25").Value = $taskTo
$ws.Range("AA25").Value = $taskStatus

# --- Row 30: first table (C:G) ---
$ws.Range("C30").Value = $taskName
$ws.Range("D30").Value = $taskDate
$ws.Range("E30").Value = $taskFrom
$ws.Range("F30").Value = $taskTo
$ws.Range("G30").Value = $taskStatus

# --- Row 30: second table (M:Q) ---
$ws.Range("M30").Value = $taskName
$ws.Range("N30").Value = $taskDate
$ws.Range("O30").Value = $taskFrom
$ws.Range("P30").Value = $taskTo
$ws.Range("Q30").Value = $taskStatus

# --- Update view: top-left visible cell & current selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = $ws.Range("J5").Row
$win.ScrollColumn = $ws.Range("J5").Column
$ws.Range("M23").Select()
